# elec category import purpose
# Clean up the "camera" row: drop a handful of stale/unused values and
# replace a few others, then move the active selection to T2 (matching
# the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Product name -> Zoom
$ws.Range("B2").Value = "Zoom"

# Other Unique code -> cleared
$ws.Range("D2").ClearContents()

# price -> 10000
$ws.Range("E2").Value = "10000"

# Special price -> 9000
$ws.Range("F2").Value = "9000"

# Offers -> cleared
$ws.Range("G2").ClearContents()

# Discount -> cleared
$ws.Range("H2").ClearContents()

# Meta keywords -> cleared
$ws.Range("J2").ClearContents()

# Meta title -> cleared
$ws.Range("K2").ClearContents()

# Ideal for -> cleared
$ws.Range("P2").ClearContents()

# Type -> cleared
$ws.Range("T2").ClearContents()

# Move the active selection / scroll focus to T2, as in the saved file.
$ws.Range("O1").Select() | Out-Null
$ws.Range("T2").Select() | Out-Null
